$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IFrames_WYSIWYGEditor")

# Update the D3 cell text (Special_Effects value for TC2)
$ws.Range("D3").Value = "Bold;Justify;right;Decrease-indent;Italic"

# Update the selection to A3:D3 with active cell A3
$ws.Activate()
$ws.Range("A3:D3").Select()
